# Update the crypto price / 1h-volume-change data that changed in this run.
# Mirrors the upstream GitHub Actions refresh of cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.183.36"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "2.475.57"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'584.07"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "'169.25"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "2.474.72"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "'4.98"
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("D13").Value = "'0.331"
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("D14").Value = "'25.62"
$ws.Range("E14").Value = "  -2.98%  "
$ws.Range("D15").Value = "2.895.68"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "67.060.55"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("D18").Value = "2.490.44"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "'11.20"
$ws.Range("E19").Value = "  -5.30%  "
$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").Value = "'354.58"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").Value = "'4.05"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'69.17"
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").Value = "'4.25"
$ws.Range("E25").Value = "  -7.37%  "
$ws.Range("E26").Value = "  -6.49%  "
$ws.Range("D27").Value = "'9.29"
$ws.Range("E27").Value = "  -6.84%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "2.591.00"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").Value = "0.0₃0911"
$ws.Range("E30").Value = "  -5.24%  "
$ws.Range("D31").Value = "'517.77"
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").Value = "'7.75"
$ws.Range("E32").Value = "  -7.11%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.25"
$ws.Range("E33").Value = "  -5.18%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'1.78"
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -7.24%  "
$ws.Range("D37").Value = "'159.26"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "'18.41"
$ws.Range("E39").Value = "  -3.62%  "
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'0.328"
$ws.Range("E42").Value = "  -6.54%  "
$ws.Range("D43").Value = "'4.81"
$ws.Range("E43").Value = "  -5.98%  "
$ws.Range("D44").Value = "'1.67"
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("D45").Value = "'2.38"
$ws.Range("E45").Value = "  -3.99%  "
$ws.Range("D46").Value = "'38.69"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("D47").Value = "'141.27"
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("D48").Value = "'3.47"
$ws.Range("E48").Value = "  -6.33%  "
$ws.Range("D49").Value = "'0.516"
$ws.Range("E49").Value = "  -6.56%  "
$ws.Range("D50").Value = "0.0₆0255"
$ws.Range("E50").Value = "  -10.63%  "
$ws.Range("D51").Value = "'1.60"
$ws.Range("E51").Value = "  -7.34%  "
